# MorganPatrone2006a - "C_Stationarygenerator_alpha_zero" : nuevos experimentos no convexos
# Regenerated numeric experiment data (x = 2.8499999999999996, y = 4.449999999999999, etc.)
# across the generator worksheets. Source values were stored as text (shared strings) in
# the workbook even when they look numeric, so we force text entry (leading apostrophe)
# and then clear the resulting "quote prefix" cell format so the cells stay plain text
# cells, matching the original authoring.

function Set-TextValue {
    param($range, [string]$text)
    $range.Value = "'" + $text
    $range.ClearFormats()
}

$wb = $excel.ActiveWorkbook

# --- Restricciones_del_lider ---
$ws2 = $wb.Worksheets.Item("Restricciones_del_lider")
$ws2.Range("A2").Value = "2.8499999999999996 - x"
Set-TextValue $ws2.Range("B2") "-3.3499999999999996"
Set-TextValue $ws2.Range("D2") "0.3"
$ws2.Range("A3").Value = "-2.8499999999999996 + x"
Set-TextValue $ws2.Range("B3") "2.3499999999999996"
Set-TextValue $ws2.Range("D3") "0.09"

# --- Restricciones_del_follower ---
$ws3 = $wb.Worksheets.Item("Restricciones_del_follower")
$ws3.Range("A2").Value = "4.449999999999999 - y"
Set-TextValue $ws3.Range("B2") "-5.449999999999999"
Set-TextValue $ws3.Range("D2") "0.19"
Set-TextValue $ws3.Range("E2") "0"
Set-TextValue $ws3.Range("F2") "0"
$ws3.Range("A3").Value = "-4.449999999999999 + y"
Set-TextValue $ws3.Range("B3") "3.4499999999999993"
Set-TextValue $ws3.Range("D3") "0.79"
Set-TextValue $ws3.Range("E3") "-0.6"
Set-TextValue $ws3.Range("F3") "-0.2"

# --- Punto_modificado ---
$ws4 = $wb.Worksheets.Item("Punto_modificado")
Set-TextValue $ws4.Range("A2") "2.8499999999999996"
Set-TextValue $ws4.Range("B2") "4.449999999999999"

# --- Vector_bf ---
# NOTE: worksheet names "Vector_bf" and "Vector_BF" differ only by case, and
# Worksheets.Item(name) resolves case-insensitively, so it would return the
# same ("Vector_bf") sheet for both names. Use positional indices instead
# (workbook order: ... 4 Punto_modificado, 5 Vector_bf, 6 Vector_BF, 7 Vector_Alpha).
$ws5 = $wb.Worksheets.Item(5)
Set-TextValue $ws5.Range("A2") "-3.4499999999999997"

# --- Vector_BF ---
$ws6 = $wb.Worksheets.Item(6)
Set-TextValue $ws6.Range("A2") "1.21"
Set-TextValue $ws6.Range("A3") "1.6"

$wb.Save()
